$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows to append starting at row 2 (columns A-K).
# L (Park Name) and M (PERMITAUTHORIZATIONNUMBER) are left blank.
$rows = @(
  @("2023/24","Yes","Registered Trapline","TR0330T006","December","Weasel","330",1,0,3,"No"),
  @("2023/24","Yes","Registered Trapline","TR0330T006","December","Flying squirrel ","330",0,0,1,"No"),
  @("2023/24","Yes","Registered Trapline","TR0330T006","December","Marten","330",0,0,1,"No"),
  @("2023/24","Yes","Registered Trapline","TR0330T006","January","Squirrel","330",0,0,1,"No"),
  @("2023/24","Yes","Registered Trapline","TR0330T006","January","Weasel","330",0,0,4,"No"),
  @("2023/24","Yes","Registered Trapline","TR0330T006","January","Fisher","330",0,1,1,"No"),
  @("2023/24","Yes","Registered Trapline","TR0330T006","February","Weasel","330",0,0,2,"No"),
  @("2023/24","Yes","Registered Trapline","TR0330T006","February","Squirrel","330",0,0,1,"No")
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = [int]$row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
    $ws.Cells.Item($r, 10).Value = $row[9]
    $ws.Cells.Item($r, 11).Value = $row[10]
    $r = $r + 1
}

# Apply wrap-text alignment to the whole new block (A2:M9), matching the
# new cell style added to the stylesheet; this also extends the sheet
# dimension to A1:M9 and touches the otherwise-blank L/M columns.
$ws.Range("A2:M9").WrapText = $true
